$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading tab character from B11 (was "\tඕ", should become "ඕ")
$ws.Range("B11").Value = "ඕ"

# Match the author's final selection/scroll state recorded in the sheet view
$ws.Range("G8").Select() | Out-Null
